$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (spreadsheet col E) value edits - imputed / removed values
$ws.Range("E3").Value = -5.7
$ws.Range("E5").ClearContents()
$ws.Range("E21").Value = -8.699999999999999
$ws.Range("E23").ClearContents()

# Remove the "RM 232" row (row 26) entirely
$ws.Rows(26).Delete()

# Remove the "SC 92" row (originally row 28, now row 27 after the previous delete)
$ws.Rows(27).Delete()

# The "SC 193" row (originally row 34) has shifted up to row 32 after the two
# row deletions above; its D-column (spreadsheet col E) value is restored.
$ws.Range("E32").Value = -6.4
